# edit.ps1 -- "imagenes mas grandes powerpoint"
#
# Slide 1: reposition the "Subtitulo 2" placeholder (Candela/Paula/Tino/
# Alejandro/Carles) and bump its run font sizes to 22pt.
#
# Slides 2,3,4,5,7,9,11: enlarge the "Marcador de contenido 3" picture
# placeholders (bigger images).

$p = $ppt.ActivePresentation

# ---- Slide 1 ----
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item("Subtítulo 2")
$shp1.Left = 506.1817626953125
$shp1.Top = 383.08087158203125
# width/height of this textbox are unchanged by the edit

$tr1 = $shp1.TextFrame.TextRange
# "Candela" / "Paula" / "Tino" / "Alejandro" / "Carles" -> paragraphs 1..5
for ($i = 1; $i -le 5; $i++) {
  $tr1.Paragraphs($i).Font.Size = 22
}

# ---- Slides 2,3,4,5,7,9,11 ----
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item("Marcador de contenido 3")
$shp.Left = 84.49512481689453
$shp.Top = 207.27284240722656
$shp.Width = 614.5795288085938
$shp.Height = 322.2954406738281

$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item("Marcador de contenido 3")
$shp.Left = 30.053308486938477
$shp.Top = 108.91913604736328
$shp.Width = 228.8557586669922
$shp.Height = 414.3536376953125

$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item("Marcador de contenido 3")
$shp.Left = 111.02102661132812
$shp.Top = 137.39779663085938
$shp.Width = 550.8882446289062
$shp.Height = 340.42041015625

$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item("Marcador de contenido 3")
$shp.Left = 53.33338928222656
$shp.Top = 111.21598815917969
$shp.Width = 201.21205139160156
$shp.Height = 419.0260009765625

$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item("Marcador de contenido 3")
$shp.Left = 53.33338928222656
$shp.Top = 138.1250457763672
$shp.Width = 200.48489379882812
$shp.Height = 372.42041015625

$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item("Marcador de contenido 3")
$shp.Left = 67.86103057861328
$shp.Top = 152.0
$shp.Width = 192.87387084960938
$shp.Height = 360.7272644042969

$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item("Marcador de contenido 3")
$shp.Left = 69.69355010986328
$shp.Top = 152.0
$shp.Width = 170.3984375
$shp.Height = 368.7840270996094
